# Replace the two tiny placeholder pictures under "Special Control Areas"
# with hyperlinks pointing at the actual SCA map images (SCA_1.jpg / SCA_2.jpg),
# exactly as the commit "yay the docx works" does:
#   <w:r><w:drawing>...</w:drawing></w:r>  ->  <w:hyperlink r:id="..."><w:r>...<w:t>URL</w:t></w:r></w:hyperlink>

$d = $word.ActiveDocument

function Replace-PictureWithHyperlink($url) {
    # Always operate on the first remaining placeholder picture - after each
    # replacement the next picture becomes index 1.
    $shape = $d.InlineShapes.Item(1)
    $shapeRange = $shape.Range
    $startPos = $shapeRange.Start

    # Insert the visible URL text right at the picture's position; this
    # replaces the inline picture in place (keeps the paragraph mark intact).
    $insertionPoint = $d.Range($startPos, $startPos)
    $insertionPoint.InsertAfter($url)

    # Turn the freshly inserted text into a real external hyperlink.
    $textRange = $d.Range($startPos, $startPos + $url.Length)
    $d.Hyperlinks.Add($textRange, $url, "", "", $url) | Out-Null
}

Replace-PictureWithHyperlink "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/SCA_1.jpg?h=100%25&w=100%25"
Replace-PictureWithHyperlink "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/SCA_2.jpg?h=100%25&w=100%25"

Write-Output ("Remaining InlineShapes: " + $d.InlineShapes.Count)
Write-Output ("Hyperlinks now: " + $d.Hyperlinks.Count)
